$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 9 -Day 6 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 171; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
